# ------------------------------------------------------------------
# Adds two new worksheets ("UsersData" and "Sheet3") with test-case
# data, clears the old tab selection from Sheet1, and leaves the new
# last sheet ("Sheet3") as the active/selected tab.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Add "UsersData" sheet right after Sheet2 -------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$usersData = $wb.Worksheets.Add($null, $afterSheet)
$usersData.Name = "UsersData"

$usersData.Range("A1").Value = "Test1234"
$usersData.Range("B1").Value = "Test1234@gmail.com"
$usersData.Range("C1").Value = "Tester"
$usersData.Range("D1").Value = "Tester"
$usersData.Range("E1").Value = "www.gmail.com"
$usersData.Range("F1").Value = "Runfast7#123"

# Hyperlinks: B1 -> mailto link, E1 -> web link (rId1 / rId2 respectively)
$usersData.Hyperlinks.Add($usersData.Range("B1"), "mailto:Test1234@gmail.com")
$usersData.Hyperlinks.Add($usersData.Range("E1"), "http://www.gmail.com")

# Re-apply the builtin Hyperlink cell style so both land on the same
# shared style index as the workbook's existing hyperlink cells.
$usersData.Range("B1").Style = "Hyperlink"
$usersData.Range("E1").Style = "Hyperlink"

# Column widths
$usersData.Columns(1).ColumnWidth = 10.666666666666666
$usersData.Columns(2).ColumnWidth = 19.333333333333332
$usersData.Columns(5).ColumnWidth = 15
$usersData.Columns(6).ColumnWidth = 21

# Selection left on the sheet (not the active tab)
[void]$usersData.Range("E9").Select()

# ---- 2. Add "Sheet3" sheet right after UsersData --------------------
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add($null, $afterSheet2)
$sheet3.Name = "Sheet3"

$sheet3.Range("A1").Value = "Test launch"
$sheet3.Range("B1").Value = 50000
$sheet3.Range("C1").Value = 200
$sheet3.Range("D1").Value = "New "
$sheet3.Range("E1").Value = "Electronic city"
$sheet3.Range("F1").Value = "immediate"
$sheet3.Range("G1").Value = "yeshwanthapur"
$sheet3.Range("H1").Value = "yeshwanthapur"
$sheet3.Range("I1").Value = 120
$sheet3.Range("J1").Value = 56
$sheet3.Range("K1").Value = 2

# Column width
$sheet3.Columns(1).ColumnWidth = 13.5

# This sheet ends up both the selected range owner (whole column B)
# and the active tab, matching the post-edit workbook state.
[void]$sheet3.Columns(2).Select()
[void]$sheet3.Select()

Write-Host "Added UsersData and Sheet3 worksheets"
